$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 14.321288329005
$ws.Range("C2").Value = 8.646604192923123
$ws.Range("D2").Value = 6.017928980072988
$ws.Range("E2").Value = 12.28302936700223
$ws.Range("G2").Value = 3.688617538353932
$ws.Range("K2").Value = 10.69774360502349
$ws.Range("L2").Value = 9.985908110676924
$ws.Range("M2").Value = 15.03365585050862
$ws.Range("N2").Value = 21.65834760351715
$ws.Range("O2").Value = 29.24596932324497
$ws.Range("B3").Value = 14.14172542177571
$ws.Range("C3").Value = 8.613804867984157
$ws.Range("D3").Value = 5.904784179569631
$ws.Range("E3").Value = 12.30714135524266
$ws.Range("G3").Value = 3.690701819405806
$ws.Range("K3").Value = 10.56478546033813
$ws.Range("L3").Value = 9.994153457459101
$ws.Range("M3").Value = 15.01340464654146
$ws.Range("N3").Value = 21.72020297583565
$ws.Range("O3").Value = 29.30031046805781
$ws.Range("B4").Value = 14.03354239880076
$ws.Range("C4").Value = 8.593243794133668
$ws.Range("D4").Value = 5.835975781286467
$ws.Range("E4").Value = 12.32341490554244
$ws.Range("G4").Value = 3.692050178124753
$ws.Range("K4").Value = 10.48454265620233
$ws.Range("L4").Value = 10.0005542474173
$ws.Range("M4").Value = 15.00337399245945
$ws.Range("N4").Value = 21.76001819161226
$ws.Range("O4").Value = 29.33906286051045
$ws.Range("B5").Value = 13.99002994537358
$ws.Range("C5").Value = 8.584760260000234
$ws.Range("D5").Value = 5.808144992623252
$ws.Range("E5").Value = 12.33041613817084
$ws.Range("G5").Value = 3.692616950371487
$ws.Range("K5").Value = 10.45223176866231
$ws.Range("L5").Value = 10.00349957539474
$ws.Range("M5").Value = 14.99989451839706
$ws.Range("N5").Value = 21.77670603330481
$ws.Range("O5").Value = 29.35620763281619
$ws.Range("B6").Value = 13.98284086269947
$ws.Range("C6").Value = 8.583345291741704
$ws.Range("D6").Value = 5.803537590960631
$ws.Range("E6").Value = 12.33160102117061
$ws.Range("G6").Value = 3.692712109265847
$ws.Range("K6").Value = 10.44689114958459
$ws.Range("L6").Value = 10.00400900872708
$ws.Range("M6").Value = 14.9993535822966
$ws.Range("N6").Value = 21.77950502511382
$ws.Range("O6").Value = 29.35913616509666
$ws.Range("B7").Value = 14.03295318573622
$ws.Range("C7").Value = 8.593129804389378
$ws.Range("D7").Value = 5.835599542245977
$ws.Range("E7").Value = 12.3235078295886
$ws.Range("G7").Value = 3.692057751672873
$ws.Range("K7").Value = 10.48410527648576
$ws.Range("L7").Value = 10.00059260424243
$ws.Range("M7").Value = 15.00332460030397
$ws.Range("N7").Value = 21.76024137416728
$ws.Range("O7").Value = 29.33928860564346
$ws.Range("B8").Value = 14.2589774233146
$ws.Range("C8").Value = 8.635383455422406
$ws.Range("D8").Value = 5.978804456809974
$ws.Range("E8").Value = 12.29103854565915
$ws.Range("G8").Value = 3.68932199312195
$ws.Range("K8").Value = 10.6516333848984
$ws.Range("L8").Value = 9.988473681974275
$ws.Range("M8").Value = 15.02617666598861
$ws.Range("N8").Value = 21.6792950474263
$ws.Range("O8").Value = 29.26358719510085
$ws.Range("B9").Value = 14.71627766173217
$ws.Range("C9").Value = 8.714832410458193
$ws.Range("D9").Value = 6.263043687533814
$ws.Range("E9").Value = 12.2390057770985
$ws.Range("G9").Value = 3.684499010433821
$ws.Range("K9").Value = 10.98951545265167
$ws.Range("L9").Value = 9.975303426027835
$ws.Range("M9").Value = 15.08988575344905
$ws.Range("N9").Value = 21.53506985348728
$ws.Range("O9").Value = 29.15794586014215
$ws.Range("B10").Value = 15.05765806737127
$ws.Range("C10").Value = 8.771045942057173
$ws.Range("D10").Value = 6.471543360563707
$ws.Range("E10").Value = 12.20785380691189
$ws.Range("G10").Value = 3.681282393769993
$ws.Range("K10").Value = 11.24117165185824
$ws.Range("L10").Value = 9.972053694644266
$ws.Range("M10").Value = 15.14795629395698
$ws.Range("N10").Value = 21.43787419598732
$ws.Range("O10").Value = 29.10651057287802
$ws.Range("B11").Value = 15.21346357127441
$ws.Range("C11").Value = 8.796133784258084
$ws.Range("D11").Value = 6.565837737080672
$ws.Range("E11").Value = 12.19521425894816
$ws.Range("G11").Value = 3.67988929807007
$ws.Range("K11").Value = 11.35591475516005
$ws.Range("L11").Value = 9.971962523924059
$ws.Range("M11").Value = 15.17675780925302
$ws.Range("N11").Value = 21.3955441008075
$ws.Range("O11").Value = 29.08880843127054
$ws.Range("B12").Value = 15.27248134409228
$ws.Range("C12").Value = 8.805562815809518
$ws.Range("D12").Value = 6.601426821641369
$ws.Range("E12").Value = 12.19064786685181
$ws.Range("G12").Value = 3.679371801591695
$ws.Range("K12").Value = 11.39936321532431
$ws.Range("L12").Value = 9.972126636652954
$ws.Range("M12").Value = 15.18800128142809
$ws.Range("N12").Value = 21.37978458362802
$ws.Range("O12").Value = 29.0829246824765
$ws.Range("B13").Value = 15.25977090142814
$ws.Range("C13").Value = 8.803535296719858
$ws.Range("D13").Value = 6.593767924041327
$ws.Range("E13").Value = 12.19162154486428
$ws.Range("G13").Value = 3.67948280804176
$ws.Range("K13").Value = 11.39000654229383
$ws.Range("L13").Value = 9.972082471653458
$ws.Range("M13").Value = 15.1855649110049
$ws.Range("N13").Value = 21.38316668709111
$ws.Range("O13").Value = 29.08415539189593
$ws.Range("B14").Value = 15.21831896076953
$ws.Range("C14").Value = 8.796910952237113
$ws.Range("D14").Value = 6.568768254535967
$ws.Range("E14").Value = 12.19483417326654
$ws.Range("G14").Value = 3.679846522420044
$ws.Range("K14").Value = 11.35948955824226
$ws.Range("L14").Value = 9.971972049627386
$ws.Range("M14").Value = 15.17767609875416
$ws.Range("N14").Value = 21.39424215116263
$ws.Range("O14").Value = 29.08830794203856
$ws.Range("B15").Value = 15.19292911090809
$ws.Range("C15").Value = 8.792844031648036
$ws.Range("D15").Value = 6.55343873617233
$ws.Range("E15").Value = 12.19683063265024
$ws.Range("G15").Value = 3.68007061382165
$ws.Range("K15").Value = 11.34079555334211
$ws.Range("L15").Value = 9.971930255378842
$ws.Range("M15").Value = 15.17288767242358
$ws.Range("N15").Value = 21.40106132011308
$ws.Range("O15").Value = 29.09095825624978
$ws.Range("B16").Value = 15.0474812821924
$ws.Range("C16").Value = 8.769396517397599
$ws.Range("D16").Value = 6.465366526032754
$ws.Range("E16").Value = 12.20871062483808
$ws.Range("G16").Value = 3.681374843335042
$ws.Range("K16").Value = 11.23367476998757
$ws.Range("L16").Value = 9.972087496299126
$ws.Range("M16").Value = 15.14612153435244
$ws.Range("N16").Value = 21.4406784024521
$ws.Range("O16").Value = 29.10778212491207
$ws.Range("B17").Value = 14.95834434049032
$ws.Range("C17").Value = 8.754887134210325
$ws.Range("D17").Value = 6.411167323720704
$ws.Range("E17").Value = 12.21639068420917
$ws.Range("G17").Value = 3.682192879544376
$ws.Range("K17").Value = 11.16799828430762
$ws.Range("L17").Value = 9.972538671061649
$ws.Range("M17").Value = 15.13030806613655
$ws.Range("N17").Value = 21.46546416679442
$ws.Range("O17").Value = 29.11956242280408
$ws.Range("B18").Value = 14.9071263874717
$ws.Range("C18").Value = 8.746496262104282
$ws.Range("D18").Value = 6.379943501379794
$ws.Range("E18").Value = 12.22095223642071
$ws.Range("G18").Value = 3.68266999904167
$ws.Range("K18").Value = 11.13025005735853
$ws.Range("L18").Value = 9.972928794318888
$ws.Range("M18").Value = 15.12143743966856
$ws.Range("N18").Value = 21.47989774108268
$ws.Range("O18").Value = 29.12687425896868
$ws.Range("B19").Value = 14.88979538842746
$ws.Range("C19").Value = 8.743647498018028
$ws.Range("D19").Value = 6.369364293342067
$ws.Range("E19").Value = 12.22252147296783
$ws.Range("G19").Value = 3.682832679749455
$ws.Range("K19").Value = 11.11747505647762
$ws.Range("L19").Value = 9.973083341549867
$ws.Range("M19").Value = 15.11847279496381
$ws.Range("N19").Value = 21.48481521206596
$ws.Range("O19").Value = 29.12944197763861
$ws.Range("B20").Value = 14.96782822231719
$ws.Range("C20").Value = 8.756436399181936
$ws.Range("D20").Value = 6.416942368976727
$ws.Range("E20").Value = 12.21555820840321
$ws.Range("G20").Value = 3.682105114789544
$ws.Range("K20").Value = 11.17498714188652
$ws.Range("L20").Value = 9.972477129916523
$ws.Range("M20").Value = 15.13196820669141
$ws.Range("N20").Value = 21.46280732121459
$ws.Range("O20").Value = 29.11825289833195
$ws.Range("B21").Value = 15.23049436051414
$ws.Range("C21").Value = 8.798858627561112
$ws.Range("D21").Value = 6.576114767960386
$ws.Range("E21").Value = 12.19388458019485
$ws.Range("G21").Value = 3.67973941864324
$ws.Range("K21").Value = 11.36845351311137
$ws.Range("L21").Value = 9.971999099304425
$ws.Range("M21").Value = 15.17998413848481
$ws.Range("N21").Value = 21.39098170021766
$ws.Range("O21").Value = 29.08706598863492
$ws.Range("B22").Value = 15.40223959107221
$ws.Range("C22").Value = 8.826168319134666
$ws.Range("D22").Value = 6.679437970094476
$ws.Range("E22").Value = 12.1810013747547
$ws.Range("G22").Value = 3.67825178926298
$ws.Range("K22").Value = 11.49486383186476
$ws.Range("L22").Value = 9.972844049100635
$ws.Range("M22").Value = 15.21332638695961
$ws.Range("N22").Value = 21.34561262792341
$ws.Range("O22").Value = 29.07146125555099
$ws.Range("B23").Value = 15.31058689219654
$ws.Range("C23").Value = 8.811631184519825
$ws.Range("D23").Value = 6.624369196491586
$ws.Range("E23").Value = 12.18776021173604
$ws.Range("G23").Value = 3.679040429720631
$ws.Range("K23").Value = 11.42741213724521
$ws.Range("L23").Value = 9.972287478323393
$ws.Range("M23").Value = 15.19535363116803
$ws.Range("N23").Value = 21.36968335679683
$ws.Range("O23").Value = 29.07935251682654
$ws.Range("B24").Value = 14.96354046766255
$ws.Range("C24").Value = 8.75573612973885
$ws.Range("D24").Value = 6.414331667355667
$ws.Range("E24").Value = 12.21593411534908
$ws.Range("G24").Value = 3.682144771989049
$ws.Range("K24").Value = 11.17182744496065
$ws.Range("L24").Value = 9.972504545376202
$ws.Range("M24").Value = 15.13121696908757
$ws.Range("N24").Value = 21.46400790813229
$ws.Range("O24").Value = 29.11884325482908
$ws.Range("B25").Value = 14.59138765933582
$ws.Range("C25").Value = 8.693712555539388
$ws.Range("D25").Value = 6.186037484577279
$ws.Range("E25").Value = 12.25183786855137
$ws.Range("G25").Value = 3.685746110127259
$ws.Range("K25").Value = 10.89734424562349
$ws.Range("L25").Value = 9.977735151561216
$ws.Range("M25").Value = 15.07065275427342
$ws.Range("N25").Value = 21.57254140668876
$ws.Range("O25").Value = 29.09095825624978
